{"js": "// Update the date heading and the 25 \"two-digit \u00f7 one-digit\" answer cells\n// in the single table, in document order, per the target diff.\n\n// 1) Update the title paragraph (first paragraph in the body).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items[0].insertText(\"2024-07-07 Sunday\", Word.InsertLocation.replace);\n\n// 2) Update the table cells. The table has 20 rows (5 cols each); only\n// every 4th row (0, 4, 8, 12, 16) holds visible answers, the rest are\n// blank spacer rows, matching the source document's layout.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst newValues = [\n  { row: 0, col: 0, text: \"20\u00f77=2, 6\" },\n  { row: 0, col: 1, text: \"55\u00f73=18, 1\" },\n  { row: 0, col: 2, text: \"33\u00f76=5, 3\" },\n  { row: 0, col: 3, text: \"29\u00f74=7, 1\" },\n  { row: 0, col: 4, text: \"66\u00f77=9, 3\" },\n\n  { row: 4, col: 0, text: \"69\u00f78=8, 5\" },\n  { row: 4, col: 1, text: \"11\u00f76=1, 5\" },\n  { row: 4, col: 2, text: \"64\u00f79=7, 1\" },\n  { row: 4, col: 3, text: \"79\u00f78=9, 7\" },\n  { row: 4, col: 4, text: \"30\u00f73=10, 0\" },\n\n  { row: 8, col: 0, text: \"34\u00f79=3, 7\" },\n  { row: 8, col: 1, text: \"66\u00f73=22, 0\" },\n  { row: 8, col: 2, text: \"32\u00f77=4, 4\" },\n  { row: 8, col: 3, text: \"94\u00f77=13, 3\" },\n  { row: 8, col: 4, text: \"98\u00f79=10, 8\" },\n\n  { row: 12, col: 0, text: \"83\u00f72=41, 1\" },\n  { row: 12, col: 1, text: \"40\u00f79=4, 4\" },\n  { row: 12, col: 2, text: \"38\u00f72=19, 0\" },\n  { row: 12, col: 3, text: \"48\u00f77=6, 6\" },\n  { row: 12, col: 4, text: \"54\u00f73=18, 0\" },\n\n  { row: 16, col: 0, text: \"27\u00f75=5, 2\" },\n  { row: 16, col: 1, text: \"93\u00f76=15, 3\" },\n  { row: 16, col: 2, text: \"58\u00f75=11, 3\" },\n  { row: 16, col: 3, text: \"70\u00f72=35, 0\" },\n  { row: 16, col: 4, text: \"90\u00f76=15, 0\" },\n];\n\nfor (const { row, col, text } of newValues) {\n  table.getCell(row, col).value = text;\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and the 25 \"two-digit \u00f7 one-digit\" answer cells\n# in the single table, in document order, per the target diff.\n\n$d = $word.ActiveDocument\n\n# 1) Update the title paragraph (first paragraph in the body).\n$d.Paragraphs.Item(1).Range.Text = \"2024-07-07 Sunday\"\n\n# 2) Update the table cells. The table has 20 rows (5 cols each); only\n# every 4th row (1, 5, 9, 13, 17 in 1-based COM indexing) holds visible\n# answers, the rest are blank spacer rows, matching the source document.\n$t = $d.Tables.Item(1)\n\n$t.Cell(1, 1).Range.Text = \"20\u00f77=2, 6\"\n$t.Cell(1, 2).Range.Text = \"55\u00f73=18, 1\"\n$t.Cell(1, 3).Range.Text = \"33\u00f76=5, 3\"\n$t.Cell(1, 4).Range.Text = \"29\u00f74=7, 1\"\n$t.Cell(1, 5).Range.Text = \"66\u00f77=9, 3\"\n\n$t.Cell(5, 1).Range.Text = \"69\u00f78=8, 5\"\n$t.Cell(5, 2).Range.Text = \"11\u00f76=1, 5\"\n$t.Cell(5, 3).Range.Text = \"64\u00f79=7, 1\"\n$t.Cell(5, 4).Range.Text = \"79\u00f78=9, 7\"\n$t.Cell(5, 5).Range.Text = \"30\u00f73=10, 0\"\n\n$t.Cell(9, 1).Range.Text = \"34\u00f79=3, 7\"\n$t.Cell(9, 2).Range.Text = \"66\u00f73=22, 0\"\n$t.Cell(9, 3).Range.Text = \"32\u00f77=4, 4\"\n$t.Cell(9, 4).Range.Text = \"94\u00f77=13, 3\"\n$t.Cell(9, 5).Range.Text = \"98\u00f79=10, 8\"\n\n$t.Cell(13, 1).Range.Text = \"83\u00f72=41, 1\"\n$t.Cell(13, 2).Range.Text = \"40\u00f79=4, 4\"\n$t.Cell(13, 3).Range.Text = \"38\u00f72=19, 0\"\n$t.Cell(13, 4).Range.Text = \"48\u00f77=6, 6\"\n$t.Cell(13, 5).Range.Text = \"54\u00f73=18, 0\"\n\n$t.Cell(17, 1).Range.Text = \"27\u00f75=5, 2\"\n$t.Cell(17, 2).Range.Text = \"93\u00f76=15, 3\"\n$t.Cell(17, 3).Range.Text = \"58\u00f75=11, 3\"\n$t.Cell(17, 4).Range.Text = \"70\u00f72=35, 0\"\n$t.Cell(17, 5).Range.Text = \"90\u00f76=15, 0\"\n"}
